$wb = $excel.ActiveWorkbook

# Rename sheet "Wong3" (sheetId 22) to "Euclid"
$ws = $wb.Worksheets.Item("Wong3")
$ws.Name = "Euclid"

# Update SBFL:RANK (col C) and SBFL:EXAM (col D) values on the renamed sheet
$updates = @(
    @{ Row = 3;  C = 2092; D = "90.25021570319241" },
    @{ Row = 5;  C = 2179; D = "94.00345125107852" },
    @{ Row = 7;  C = 1898; D = "81.88093183779121" },
    @{ Row = 10; C = 2179; D = "94.00345125107852" },
    @{ Row = 13; C = 1716; D = "74.02933563416738" },
    @{ Row = 15; C = 2146; D = "92.57981018119068" },
    @{ Row = 18; C = 1576; D = "67.69759450171821" },
    @{ Row = 19; C = 2073; D = "89.43054357204487" },
    @{ Row = 20; C = 1929; D = "83.14655172413794" },
    @{ Row = 25; C = 2141; D = "92.36410698878343" },
    @{ Row = 28; C = 2132; D = "91.97584124245039" },
    @{ Row = 29; C = 1572; D = "68.02250108178278" },
    @{ Row = 30; C = 1973; D = "84.75085910652921" },
    @{ Row = 33; C = 18;   D = "0.7712082262210797" },
    @{ Row = 34; C = 2;    D = "0.0859106529209622" },
    @{ Row = 35; C = 1914; D = "82.21649484536083" },
    @{ Row = 36; C = 2164; D = "93.35634167385678" },
    @{ Row = 42; C = 2148; D = "92.66609145815357" },
    @{ Row = 43; C = 1793; D = "77.35116479723901" },
    @{ Row = 46; C = 2185; D = "94.26229508196722" },
    @{ Row = 50; C = 1789; D = "77.17860224331319" },
    @{ Row = 57; C = 2072; D = "89.38740293356342" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 4).Value = $u.D
}
